$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Replace run text: drop the trailing-space run and update the ID text.
$d.Content.Find.Execute("**ID__AFFARS_5316_topic_11__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5316_301_3__ID**", 2) | Out-Null

# Update left indent: 120 twips -> 225 twips (6pt -> 11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Add a paragraph border (space-only, no line) around the first paragraph only.
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

Write-Output "Done: [$($p1.Range.Text)]"
